# Insert a new data row at row 254 (pushing existing rows 254:301 down to
# 255:302) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 254..301 down to 255..302, leaving a blank row 254 to fill in.
$ws.Rows.Item(254).Insert()

# Columns that stay constant for every data row in this sheet.
$ws.Cells.Item(254, 1).Value = 11
$ws.Cells.Item(254, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(254, 3).Value = "Bíobío"
$ws.Cells.Item(254, 5).Value = 8
$ws.Cells.Item(254, 6).Value = 100112045
$ws.Cells.Item(254, 7).Value = "Zapallo"
$ws.Cells.Item(254, 17).Value = 1
$ws.Cells.Item(254, 18).Value = "Hortaliza"

# New record's own data.
$ws.Cells.Item(254, 4).Value = 44889
$ws.Cells.Item(254, 8).Value = "Paine"
$ws.Cells.Item(254, 9).Value = "2a nueva(o)"
$ws.Cells.Item(254, 10).Value = 1000
$ws.Cells.Item(254, 11).Value = 450
$ws.Cells.Item(254, 12).Value = 500
$ws.Cells.Item(254, 13).Value = 475
$ws.Cells.Item(254, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(254, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(254, 16).Value = 475
